# Move the "Dla wydatkow indywidualnych" summary block (columns L:N, rows 1-3)
# down to columns G:I, rows 7-9, on every monthly-category worksheet (sheets
# 2..8 - the "Ogolne" sheet, index 1, is untouched since it never had this
# block). The block keeps its 3x3 shape (label+2 blanks / 3 headers / 3
# numbers) but is re-homed under the "Dla sum miesiecznych" block instead of
# sitting beside it, and columns L/M no longer need a custom width.

$wb = $excel.ActiveWorkbook

for ($idx = 2; $idx -le $wb.Worksheets.Count; $idx++) {
    $ws = $wb.Worksheets.Item($idx)

    # Nothing to do if this sheet never had the L:N block.
    if ($ws.Range("L1").Value2 -eq $null -and $ws.Range("L1").MergeCells -ne $true) {
        continue
    }

    # --- 1. Snapshot the old values (before anything is touched) ---
    $g7v = $ws.Range("L1").Value2
    $h7v = $ws.Range("M1").Value2
    $i7v = $ws.Range("N1").Value2
    $g8v = $ws.Range("L2").Value2
    $h8v = $ws.Range("M2").Value2
    $i8v = $ws.Range("N2").Value2
    $g9v = $ws.Range("L3").Value2
    $h9v = $ws.Range("M3").Value2
    $i9v = $ws.Range("N3").Value2

    # --- 2. Stamp the destination rows with the same look as the existing
    #        "Dla sum miesiecznych" block (G1:I1 / G2:I2 / G3:I3 use style
    #        index 2 - thin border, centered) ---
    $ws.Range("G1:I1").Copy()
    $ws.Range("G7:I7").PasteSpecial(-4122)
    $ws.Range("G1:I1").Copy()
    $ws.Range("G8:I8").PasteSpecial(-4122)
    $ws.Range("G1:I1").Copy()
    $ws.Range("G9:I9").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- 3. Merge G7:I7 like L1:N1 used to be merged. Merging re-derives
    #        per-cell borders for the (now hidden) inner cells, so re-apply
    #        the format template once more afterwards to keep plain style 2
    #        on all three cells (matching how G1:I1 looks despite also
    #        being merged). ---
    $ws.Range("G7:I7").Merge()
    $ws.Range("G1:I1").Copy()
    $ws.Range("G7:I7").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # --- 4. Write the values into their new home ---
    if ($g7v -ne $null) { $ws.Range("G7").Value2 = $g7v }
    if ($h7v -ne $null) { $ws.Range("H7").Value2 = $h7v }
    if ($i7v -ne $null) { $ws.Range("I7").Value2 = $i7v }
    if ($g8v -ne $null) { $ws.Range("G8").Value2 = $g8v }
    if ($h8v -ne $null) { $ws.Range("H8").Value2 = $h8v }
    if ($i8v -ne $null) { $ws.Range("I8").Value2 = $i8v }
    if ($g9v -ne $null) { $ws.Range("G9").Value2 = $g9v }
    if ($h9v -ne $null) { $ws.Range("H9").Value2 = $h9v }
    if ($i9v -ne $null) { $ws.Range("I9").Value2 = $i9v }

    # --- 5. Remove the old L1:N3 block entirely (contents + formatting),
    #        then drop columns L and M's custom-width definitions. Deleting
    #        starting at column J (which has no custom width of its own)
    #        instead of directly at L avoids leaving a stray/invalid <col>
    #        entry behind. ---
    $ws.Range("L1:N1").UnMerge()
    $ws.Range("L1:N3").Clear()
    $ws.Range("J1:M1").EntireColumn.Delete()
}

Write-Output "Moved individual-expense summary blocks on all category sheets."
